# Auto-generated Excel COM-interop script
# Applies updated values to pl_mw.xlsx data table (rows 2-25, columns B,C,E,F,G,H,I,J,K,L,N)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8026254900991603
$ws.Range("C2").Value = 0.06972833565704661
$ws.Range("E2").Value = 0.1749855810264549
$ws.Range("F2").Value = 2.299579057228655
$ws.Range("G2").Value = 1.202446591375022
$ws.Range("H2").Value = 1.173121123765128
$ws.Range("I2").Value = 1.179346095783742
$ws.Range("J2").Value = 0.07156967794969304
$ws.Range("K2").Value = 0.4774595367006782
$ws.Range("L2").Value = 0.4501545214245368
$ws.Range("N2").Value = 2.116977087754879
$ws.Range("B3").Value = 0.7629333717115117
$ws.Range("C3").Value = 0.06753713348266643
$ws.Range("E3").Value = 0.1739920969943221
$ws.Range("F3").Value = 2.29615214440652
$ws.Range("G3").Value = 1.205796194450144
$ws.Range("H3").Value = 1.179304178624633
$ws.Range("I3").Value = 1.185879927176849
$ws.Range("J3").Value = 0.0711796722019109
$ws.Range("K3").Value = 0.4410672433099592
$ws.Range("L3").Value = 0.4415976814238007
$ws.Range("N3").Value = 2.137620415348501
$ws.Range("B4").Value = 0.7388949428227249
$ws.Range("C4").Value = 0.06616989730799361
$ws.Range("E4").Value = 0.1734417361573009
$ws.Range("F4").Value = 2.295131348574216
$ws.Range("G4").Value = 1.208498619429449
$ws.Range("H4").Value = 1.183559340922798
$ws.Range("I4").Value = 1.190403626308306
$ws.Range("J4").Value = 0.07094042061022066
$ws.Range("K4").Value = 0.4188926129581034
$ws.Range("L4").Value = 0.4365346361910838
$ws.Range("N4").Value = 2.150940131486735
$ws.Range("B5").Value = 0.7291833438817719
$ws.Range("C5").Value = 0.06560722890045412
$ws.Range("E5").Value = 0.1732325115777051
$ws.Range("F5").Value = 2.294988044170225
$ws.Range("G5").Value = 1.209762212418525
$ws.Range("H5").Value = 1.18540878885527
$ws.Range("I5").Value = 1.192375835219174
$ws.Range("J5").Value = 0.07084298698290858
$ws.Range("K5").Value = 0.4098995030083472
$ws.Range("L5").Value = 0.4345195896189153
$ws.Range("N5").Value = 2.156530116192238
$ws.Range("B6").Value = 0.7275758467594358
$ws.Range("C6").Value = 0.06551346496945598
$ws.Range("E6").Value = 0.1731986806948171
$ws.Range("F6").Value = 2.294980725439714
$ws.Range("G6").Value = 1.209981834214602
$ws.Range("H6").Value = 1.185722862633185
$ws.Range("I6").Value = 1.192711097375231
$ws.Range("J6").Value = 0.07082681231178078
$ws.Range("K6").Value = 0.4084088264898185
$ws.Range("L6").Value = 0.4341879083973197
$ws.Range("N6").Value = 2.15746811876436
$ws.Range("B7").Value = 0.7387636268600488
$ws.Range("C7").Value = 0.0661623312937536
$ws.Range("E7").Value = 0.1734388534643649
$ws.Range("F7").Value = 2.295128311510013
$ws.Range("G7").Value = 1.208515003471177
$ws.Range("H7").Value = 1.183583815728596
$ws.Range("I7").Value = 1.190429702787718
$ws.Range("J7").Value = 0.07093910631740918
$ws.Range("K7").Value = 0.4187711531736511
$ws.Range("L7").Value = 0.4365072652205697
$ws.Range("N7").Value = 2.151014863678771
$ws.Range("B8").Value = 0.7888709419624718
$ws.Range("C8").Value = 0.06897732758496034
$ws.Range("E8").Value = 0.1746306809380833
$ws.Range("F8").Value = 2.29817276081917
$ws.Range("G8").Value = 1.203467476992017
$ws.Range("H8").Value = 1.175157876681894
$ws.Range("I8").Value = 1.181492740990905
$ws.Range("J8").Value = 0.07143516532537575
$ws.Range("K8").Value = 0.4648763308144055
$ws.Range("L8").Value = 0.4471646112379233
$ws.Range("N8").Value = 2.123961113075465
$ws.Range("B9").Value = 0.889750369075216
$ws.Range("C9").Value = 0.07432540325916648
$ws.Range("E9").Value = 0.1774389544640798
$ws.Range("F9").Value = 2.312730247377189
$ws.Range("G9").Value = 1.198696385671639
$ws.Range("H9").Value = 1.1622713970239
$ws.Range("I9").Value = 1.168027658917715
$ws.Range("J9").Value = 0.07240925107964458
$ws.Range("K9").Value = 0.5566282253591055
$ws.Range("L9").Value = 0.4695716142282436
$ws.Range("N9").Value = 2.076023046130285
$ws.Range("B10").Value = 0.9654437686828885
$ws.Range("C10").Value = 0.07815149388008535
$ws.Range("E10").Value = 0.1797869028766961
$ws.Range("F10").Value = 2.328653774860499
$ws.Range("G10").Value = 1.198322961489652
$ws.Range("H10").Value = 1.155017110328785
$ws.Range("I10").Value = 1.160608907477453
$ws.Range("J10").Value = 0.07312527160646809
$ws.Range("K10").Value = 0.6248459731415323
$ws.Range("L10").Value = 0.486946606634902
$ws.Range("N10").Value = 2.043918741923797
$ws.Range("B11").Value = 1.00021754822356
$ws.Range("C11").Value = 0.07987008628116143
$ws.Range("E11").Value = 0.1809163642912743
$ws.Range("F11").Value = 2.337031887592573
$ws.Range("G11").Value = 1.1988346242034
$ws.Range("H11").Value = 1.152196858942446
$ws.Range("I11").Value = 1.157770959278025
$ws.Range("J11").Value = 0.07345099567083579
$ws.Range("K11").Value = 0.6560538663777891
$ws.Range("L11").Value = 0.4950478017074289
$ws.Range("N11").Value = 2.029989810120096
$ws.Range("B12").Value = 1.013433944724909
$ws.Range("C12").Value = 0.08051774617423035
$ws.Range("E12").Value = 0.1813528369099195
$ws.Range("F12").Value = 2.340367406906978
$ws.Range("G12").Value = 1.199126475248718
$ws.Range("H12").Value = 1.151197830646353
$ws.Range("I12").Value = 1.156773484456025
$ws.Range("J12").Value = 0.07357433013837777
$ws.Range("K12").Value = 0.667896400453543
$ws.Range("L12").Value = 0.4981437177685564
$ws.Range("N12").Value = 2.024812434320465
$ws.Range("B13").Value = 1.010585418346352
$ws.Range("C13").Value = 0.08037840012141118
$ws.Range("E13").Value = 0.1812584455377753
$ws.Range("F13").Value = 2.339641801160411
$ws.Range("G13").Value = 1.199059255619829
$ws.Range("H13").Value = 1.151409923926778
$ws.Range("I13").Value = 1.156984875376303
$ws.Range("J13").Value = 0.07354776847455113
$ws.Range("K13").Value = 0.6653448034277005
$ws.Range("L13").Value = 0.4974757070611133
$ws.Range("N13").Value = 2.025923150226305
$ws.Range("B14").Value = 1.001303904352312
$ws.Range("C14").Value = 0.07992343235389399
$ws.Range("E14").Value = 0.1809520977027752
$ws.Range("F14").Value = 2.337303039068971
$ws.Range("G14").Value = 1.198856668459726
$ws.Range("H14").Value = 1.152113286907976
$ws.Range("I14").Value = 1.157687349373319
$ws.Range("J14").Value = 0.07346114273621751
$ws.Range("K14").Value = 0.6570276645479112
$ws.Range("L14").Value = 0.4953019411259305
$ws.Range("N14").Value = 2.029561915652076
$ws.Range("B15").Value = 0.9956249828904333
$ws.Range("C15").Value = 0.07964434387722008
$ws.Range("E15").Value = 0.1807655913155948
$ws.Range("F15").Value = 2.335891688947299
$ws.Range("G15").Value = 1.198745355623728
$ws.Range("H15").Value = 1.152553093581574
$ws.Range("I15").Value = 1.158127687766267
$ws.Range("J15").Value = 0.07340808033242041
$ws.Range("K15").Value = 0.6519363915811311
$ws.Range("L15").Value = 0.4939741095218295
$ws.Range("N15").Value = 2.031803426240876
$ws.Range("B16").Value = 0.9631780169828801
$ws.Range("C16").Value = 0.0780387403416114
$ws.Range("E16").Value = 0.1797143201330833
$ws.Range("F16").Value = 2.328129054255768
$ws.Range("G16").Value = 1.198303245177371
$ws.Range("H16").Value = 1.155211070159325
$ws.Range("I16").Value = 1.160805176747211
$ws.Range("J16").Value = 0.07310398392553097
$ws.Range("K16").Value = 0.6228099549981607
$ws.Range("L16").Value = 0.4864211267868797
$ws.Range("N16").Value = 2.044842629088571
$ws.Range("B17").Value = 0.943359602967746
$ws.Range("C17").Value = 0.07704815408649779
$ws.Range("E17").Value = 0.1790850782861462
$ws.Range("F17").Value = 2.323657326144613
$ws.Range("G17").Value = 1.198206647300481
$ws.Range("H17").Value = 1.156964495438572
$ws.Range("I17").Value = 1.162585231219239
$ws.Range("J17").Value = 0.07291742382214395
$ws.Range("K17").Value = 0.6049864248234087
$ws.Range("L17").Value = 0.4818379942439179
$ws.Range("N17").Value = 2.053014816351576
$ws.Range("B18").Value = 0.9319926604254931
$ws.Range("C18").Value = 0.07647633123448827
$ws.Range("E18").Value = 0.1787289345144991
$ws.Range("F18").Value = 2.32119209645542
$ws.Range("G18").Value = 1.198215230578782
$ws.Range("H18").Value = 1.158018179440035
$ws.Range("I18").Value = 1.163659605493251
$ws.Range("J18").Value = 0.07281012027385714
$ws.Range("K18").Value = 0.5947513269505009
$ws.Range("L18").Value = 0.4792204727119724
$ws.Range("N18").Value = 2.057778816688771
$ws.Range("B19").Value = 0.9281495390928569
$ws.Range("C19").Value = 0.07628236683881084
$ws.Range("E19").Value = 0.1786093446490931
$ws.Range("F19").Value = 2.320375759993411
$ws.Range("G19").Value = 1.198229150712422
$ws.Range("H19").Value = 1.158382696556643
$ws.Range("I19").Value = 1.164032049686888
$ws.Range("J19").Value = 0.0727737896109879
$ws.Range("K19").Value = 0.5912887498664929
$ws.Range("L19").Value = 0.4783374210103659
$ws.Range("N19").Value = 2.059402742061938
$ws.Range("B20").Value = 0.94546599092007
$ws.Range("C20").Value = 0.07715381715613745
$ws.Range("E20").Value = 0.1791514643944225
$ws.Range("F20").Value = 2.324122298446781
$ws.Range("G20").Value = 1.198210290999924
$ws.Range("H20").Value = 1.156773166774244
$ws.Range("I20").Value = 1.162390511358055
$ws.Range("J20").Value = 0.07293728344809125
$ws.Range("K20").Value = 0.6068820628161404
$ws.Range("L20").Value = 0.4823239554614247
$ws.Range("N20").Value = 2.052138292816107
$ws.Range("B21").Value = 1.004028803671986
$ws.Range("C21").Value = 0.08005715232415156
$ws.Range("E21").Value = 0.1810418419045021
$ws.Range("F21").Value = 2.337985570601589
$ws.Range("G21").Value = 1.198913510205116
$ws.Range("H21").Value = 1.151904821534629
$ws.Range("I21").Value = 1.157478920753874
$ws.Range("J21").Value = 0.07348658717537404
$ws.Range("K21").Value = 0.6594699406553559
$ws.Range("L21").Value = 0.4959396658919673
$ws.Range("N21").Value = 2.028490483215212
$ws.Range("B22").Value = 1.042584342193322
$ws.Range("C22").Value = 0.08193639197531866
$ws.Range("E22").Value = 0.1823284064287911
$ws.Range("F22").Value = 2.347995482312712
$ws.Range("G22").Value = 1.199944922847564
$ws.Range("H22").Value = 1.149124887632198
$ws.Range("I22").Value = 1.15471885817837
$ws.Range("J22").Value = 0.07384552474347217
$ws.Range("K22").Value = 0.693983411077113
$ws.Range("L22").Value = 0.5050024232019439
$ws.Range("N22").Value = 2.013601914815263
$ws.Range("B23").Value = 1.02198099257356
$ws.Range("C23").Value = 0.08093507115934528
$ws.Range("E23").Value = 0.1816370856726799
$ws.Range("F23").Value = 2.342566200373341
$ws.Range("G23").Value = 1.19934208677526
$ws.Range("H23").Value = 1.150571840946924
$ws.Range("I23").Value = 1.156150787908942
$ws.Range("J23").Value = 0.07365396237393185
$ws.Range("K23").Value = 0.6755498638313782
$ws.Range("L23").Value = 0.5001505068042746
$ws.Range("N23").Value = 2.021496350360106
$ws.Range("B24").Value = 0.9445136090012909
$ws.Range("C24").Value = 0.07710605411463689
$ws.Range("E24").Value = 0.1791214337388567
$ws.Range("F24").Value = 2.323911755431993
$ws.Range("G24").Value = 1.198208443961619
$ws.Range("H24").Value = 1.156859524363512
$ws.Range("I24").Value = 1.162478385350482
$ws.Range("J24").Value = 0.07292830505833336
$ws.Range("K24").Value = 0.6060250078215006
$ws.Range("L24").Value = 0.4821041982375789
$ws.Range("N24").Value = 2.052534364450507
$ws.Range("B25").Value = 0.8621814238216245
$ws.Range("C25").Value = 0.07289684896372961
$ws.Range("E25").Value = 0.1766290645612116
$ws.Range("F25").Value = 2.307873445865226
$ws.Range("G25").Value = 1.199437554311004
$ws.Range("H25").Value = 1.165368572879345
$ws.Range("I25").Value = 1.171235768112389
$ws.Range("J25").Value = 0.07214563989277067
$ws.Range("K25").Value = 0.5316644192079991
$ws.Range("L25").Value = 0.4633491829498269
$ws.Range("N25").Value = 2.088444266246983
